# Label BOM items better.
# Relabel resistor/capacitor Package and Description columns with clearer,
# more consistent naming conventions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resistors (rows 2-11): Package "R-W4" -> "R-1/4W" ---
$ws.Range("C2:C11").Value = "R-1/4W"

# --- Ceramic capacitors (rows 12-13): Package "C-5mm" -> "C-P5mm" ---
$ws.Range("C12:C13").Value = "C-P5mm"

# --- Film capacitors (rows 14-15): Package relabeled with new "C-P5mm" prefix ---
$ws.Range("C14").Value = "C-P5mm 2.5x7.2mm"
$ws.Range("C15").Value = "C-P5mm 5x7.2 / 11x7.2 mm"

# --- Descriptions: reorder "Capacitor X THT" -> "X Capacitor THT" ---
$ws.Range("E14:E15").Value = "Film Capacitor THT"
$ws.Range("E12:E13").Value = "Ceramic Capacitor THT"

# --- Electrolytic capacitor (row 16): Package + Description relabeled ---
$ws.Range("E16").Value = "Electrolytic Capacitor THT"
$ws.Range("C16").Value = "E-P2.5mm 6.3x11.5mm"

# --- Column C width widened slightly to fit the longer labels ---
$ws.Columns.Item(3).ColumnWidth = 26

# --- Restore active selection to C16 (last cell touched in the edit) ---
$ws.Range("C16").Select()
